$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: who / which week ---
$ws.Range("C1").Value = "Jesse Hare"
$ws.Range("E1").Value = 7

# --- New task rows (3-5) ---
$ws.Range("A3").Value = "Project Build"
$ws.Range("B3").Value = "Finish work on first iteration"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 3

$ws.Range("A4").Value = "Project Build"
$ws.Range("B4").Value = "Iteration Review"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1

$ws.Range("A5").Value = "Project Build"
$ws.Range("B5").Value = "Gather any new requirements, commence work on next iteration"
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 16

# --- Cumulative total label now carries the computed total ---
$ws.Range("A14").Value = "Cumulative Total: 140"

# --- Widen column B to fit the longer task descriptions ---
$ws.Columns.Item(2).ColumnWidth = 48.9

# --- Selection moved to B12 ---
$ws.Range("B12").Select() | Out-Null
